$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.483.71'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '2.444.04'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''555.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.38%  '
$ws.Range('D6').Value = '''160.67'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '''0.500'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.06%  '
$ws.Range('D9').Value = '2.441.60'
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('D10').Value = '''0.148'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.17%  '
$ws.Range('D11').Value = '''0.165'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('D12').Value = '''0.333'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.19%  '
$ws.Range('D13').Value = '''4.77'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.02%  '
$ws.Range('D14').Value = '2.884.06'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').Value = '68.366.01'
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').Value = '''0.0000167'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.23%  '
$ws.Range('D17').Value = '''23.23'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.63%  '
$ws.Range('D18').Value = '2.448.28'
$ws.Range('E18').Value = '  -1.73%  '
$ws.Range('D19').Value = '''10.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.17%  '
$ws.Range('D20').Value = '''338.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.31%  '
$ws.Range('D21').Value = '''6.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.51%  '
$ws.Range('D22').Value = '''3.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.40%  '
$ws.Range('D23').Value = '''5.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '''1.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.98%  '
$ws.Range('D26').Value = '''66.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.40%  '
$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').Value = '''3.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.17%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.563.06'
$ws.Range('E28').Value = '  -2.05%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '''8.05'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.02%  '
$ws.Range('D31').Value = '0.0₃0808'
$ws.Range('E31').Value = '  -6.78%  '
$ws.Range('D32').Value = '''7.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.82%  '
$ws.Range('D33').Value = '''1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '''432.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('D35').Value = '''1.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.74%  '
$ws.Range('D36').Value = '''1.60'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.93%  '
$ws.Range('D37').Value = '''156.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('D38').Value = '''19.03'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = '''0.109'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.98%  '
$ws.Range('D41').Value = '''17.74'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.93%  '
$ws.Range('D42').Value = '''0.301'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.61%  '
$ws.Range('D43').Value = '''4.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.15%  '
$ws.Range('D44').Value = '''37.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('D45').Value = '''1.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.94%  '
$ws.Range('D46').Value = '''1.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('D47').Value = '''2.03'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.95%  '
$ws.Range('D48').Value = '''131.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.03%  '
$ws.Range('D49').Value = '''3.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.80%  '
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('D51').Value = '''0.479'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.87%  '
